# harmonized similar tags to be the same
# The "Tags Term Accession Number" (B13) on the SwateTemplateMetadata sheet
# used a full PURL ("http://purl.obolibrary.org/obo/NCIT_C17156"); it is
# harmonized to the short CURIE form "NCIT:C17156" used elsewhere in the
# template. The now-redundant "Tags Term Source REF" (B14) value "NCIT" is
# cleared since the source ref is now embedded in the accession itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SwateTemplateMetadata")

$ws.Range("B13").Value = "NCIT:C17156"
$ws.Range("B14").Value = ""

$ws.Activate() | Out-Null
$ws.Range("B14").Select() | Out-Null
